$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update B2 value
$ws.Range("B2").Value = 3783

# Rebuild rows 3-12 with the new labels/values (row 1 and row 2's A2 label unchanged)
$labels = @(
    "quartz monzonite",
    "syeno granite",
    "granodiorite",
    "quartz syenite",
    "monzonite",
    "quartz monzodiorite`nquartz monzogabbro",
    "quartz-rich granitoid",
    "syenite",
    "tonalite",
    "monzodiorite monzogabbro"
)

$values = @(549, 142, 92, 41, 34, 8, 5, 2, 2, 1)

for ($i = 0; $i -lt $labels.Length; $i++) {
    $row = 3 + $i
    $ws.Cells.Item($row, 1).Value = $labels[$i]
    $ws.Cells.Item($row, 2).Value = $values[$i]
}
